$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1591
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F4").Value = 365
$ws1.Range("F5").Value = 5183
$ws1.Range("F6").Value = 557
$ws1.Range("F7").Value = 10266
$ws1.Range("F8").Value = 264
$ws1.Range("F9").Value = 560
$ws1.Range("F10").Value = 109
$ws1.Range("F11").Value = 93
$ws1.Range("F12").Value = 800

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 16

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1591
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F4").Value = 365
$ws4.Range("F7").Value = 5183
$ws4.Range("F8").Value = 557
$ws4.Range("F9").Value = 16
$ws4.Range("F10").Value = 10266
$ws4.Range("F11").Value = 264
$ws4.Range("F12").Value = 560
$ws4.Range("F13").Value = 109
$ws4.Range("F16").Value = 93
$ws4.Range("F17").Value = 800
